# CS235 Assignment #4 presentation update
# - Slide 4 ("Design Patterns & Techniques"): merge the 3 runs of the
#   "Mobile friendly version of day, week, and month calendar views"
#   bullet into a single run, sized 18pt and colored dark green (003300).
# - Slide 5 ("Design Patterns & Techniques (Continued)"): give the
#   "Motion and coloring of errors", "Mobile friendly", "Calendar
#   navigation" and "Don't make 'em think!" bullets 90% within-paragraph
#   line spacing, 18pt size, and dark-green (003300) color.

$p = $ppt.ActivePresentation

# --- Slide 4 -----------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$body4 = $slide4.Shapes.Item(2)
$tr4 = $body4.TextFrame.TextRange

$calParagraph = $tr4.Paragraphs(5)
$calRange = $tr4.Characters($calParagraph.Start, $calParagraph.Length)
$calRange.Text = "Mobile friendly version of day, week, and month calendar views"

$calParagraph = $tr4.Paragraphs(5)
$calParagraph.Font.Size = 18
$calParagraph.Font.Color.RGB = 13056

# --- Slide 5 -------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$body5 = $slide5.Shapes.Item(2)
$tr5 = $body5.TextFrame.TextRange

foreach ($paraIndex in 3, 5, 7, 9) {
    $para = $tr5.Paragraphs($paraIndex)
    $para.ParagraphFormat.SpaceWithin = 0.9
    $para.Font.Size = 18
    $para.Font.Color.RGB = 13056
}
